$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "23.776.61"
$ws.Range("E2").Value = "  +1.93%  "
Set-TextValue "D3" "1.654.62"
$ws.Range("E3").Value = "  +1.87%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue "D5" "0.9995"
$ws.Range("E5").Value = "  -0.14%  "
Set-TextValue "D6" "304.11"
$ws.Range("E6").Value = "  +0.49%  "
Set-TextValue "D7" "0.3828"
$ws.Range("E7").Value = "  +2.13%  "
Set-TextValue "D8" "0.3609"
$ws.Range("E8").Value = "  -0.32%  "
Set-TextValue "D9" "51.18"
$ws.Range("E9").Value = "  -0.18%  "
Set-TextValue "D10" "1.247"
$ws.Range("E10").Value = "  +1.97%  "
Set-TextValue "D11" "0.08230"
$ws.Range("E11").Value = "  +1.10%  "
Set-TextValue "D12" "1.000"
$ws.Range("E12").Value = "  -0.09%  "
Set-TextValue "D13" "22.70"
$ws.Range("E13").Value = "  +1.85%  "
Set-TextValue "D14" "6.542"
$ws.Range("E14").Value = "  +1.11%  "
Set-TextValue "D15" "7.412"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("E16").Value = "  +0.23%  "
Set-TextValue "D17" "1.654.09"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  +3.92%  "
Set-TextValue "D19" "0.06981"
$ws.Range("E19").Value = "  +0.77%  "
Set-TextValue "D20" "6.783"
$ws.Range("E20").Value = "  +3.95%  "
$ws.Range("E21").Value = "  +1.33%  "
Set-TextValue "D22" "0.9987"
Set-TextValue "D24" "23.774.32"
$ws.Range("E24").Value = "  +1.88%  "
Set-TextValue "D25" "2.550"
$ws.Range("E25").Value = "  +3.31%  "
Set-TextValue "D26" "3.103"
$ws.Range("E26").Value = "  +0.69%  "
Set-TextValue "D27" "21.32"
$ws.Range("E27").Value = "  +1.00%  "
Set-TextValue "D28" "150.98"
$ws.Range("E28").Value = "  +0.26%  "
Set-TextValue "D29" "5.232"
$ws.Range("E29").Value = "  -0.69%  "
Set-TextValue "D30" "134.78"
$ws.Range("E30").Value = "  +1.54%  "
Set-TextValue "D31" "1.839.35"
$ws.Range("E31").Value = "  +2.20%  "
Set-TextValue "D32" "6.952"
$ws.Range("E32").Value = "  +3.19%  "
Set-TextValue "D33" "1.082"
$ws.Range("E33").Value = "  +1.92%  "
Set-TextValue "D34" "11.96"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("E35").Value = "  -2.28%  "
Set-TextValue "D36" "0.02843"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D37" "0.2521"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D38" "6.149"
$ws.Range("E38").Value = "  +2.62%  "
Set-TextValue "D39" "0.08842"
$ws.Range("E39").Value = "  +1.01%  "
Set-TextValue "D40" "0.07068"
$ws.Range("E40").Value = "  -0.47%  "
Set-TextValue "D41" "12.87"
$ws.Range("E41").Value = "  +6.75%  "
Set-TextValue "D42" "0.7090"
$ws.Range("E42").Value = "  +1.74%  "
Set-TextValue "D43" "1.340"
$ws.Range("E43").Value = "  +0.47%  "
Set-TextValue "D44" "15.84"
$ws.Range("E44").Value = "  -1.43%  "
Set-TextValue "D45" "0.6564"
$ws.Range("E45").Value = "  +1.59%  "
Set-TextValue "D46" "2.337"
$ws.Range("E46").Value = "  +3.31%  "
Set-TextValue "D47" "0.9992"
$ws.Range("E47").Value = "  -0.11%  "
Set-TextValue "D48" "3.973"
$ws.Range("E48").Value = "  +0.38%  "
Set-TextValue "D49" "0.07987"
$ws.Range("E49").Value = "  +0.28%  "
Set-TextValue "D50" "128.84"
$ws.Range("E50").Value = "  +2.52%  "
Set-TextValue "D51" "1.199"
$ws.Range("E51").Value = "  +1.30%  "
